$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("double")

# Add a value smaller than Long.MIN_VALUE to A9 (the negative of the
# "bigger than max long" double already sitting in A8), so the integration
# test can check values below Long.MIN_VALUE as well as above Long.MAX_VALUE.
$ws.Range("A9").Value = -9523372036854769700.0

# Leave the selection where it ended up after entering the new value.
$ws.Range("C6").Select() | Out-Null
